# Applies the cryptos list refresh described in the commit diff.
# Numeric-looking "Price" strings (column D) get a leading apostrophe so
# Excel stores them as literal text (matching the original inlineStr cells)
# instead of silently parsing them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.069.24'
$ws.Range("E2").Value = '  +4.32%  '
$ws.Range("D3").Value = '1.905.92'
$ws.Range("E3").Value = '  +5.30%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''251.74'
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '''0.5093'
$ws.Range("E7").Value = '  +2.49%  '
$ws.Range("D8").Value = '''45.03'
$ws.Range("E8").Value = '  +4.36%  '
$ws.Range("D9").Value = '''0.3021'
$ws.Range("E9").Value = '  +8.65%  '
$ws.Range("D10").Value = '''0.06795'
$ws.Range("E10").Value = '  +5.93%  '
$ws.Range("D11").Value = '1.907.74'
$ws.Range("E11").Value = '  +5.42%  '
$ws.Range("D12").Value = '''17.24'
$ws.Range("E12").Value = '  +2.95%  '
$ws.Range("D13").Value = '''0.07329'
$ws.Range("E13").Value = '  +3.64%  '
$ws.Range("D14").Value = '''0.6965'
$ws.Range("E14").Value = '  +7.84%  '
$ws.Range("D15").Value = '''86.52'
$ws.Range("E15").Value = '  +2.85%  '
$ws.Range("D16").Value = '''4.908'
$ws.Range("E16").Value = '  +4.64%  '
$ws.Range("D17").Value = '30.066.05'
$ws.Range("E17").Value = '  +4.33%  '
$ws.Range("D18").Value = '''0.000008155'
$ws.Range("E18").Value = '  +11.18%  '
$ws.Range("D19").Value = '''1.000'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '''13.03'
$ws.Range("E20").Value = '  +6.39%  '
$ws.Range("D21").Value = '2.153.09'
$ws.Range("E21").Value = '  +5.10%  '
$ws.Range("D22").Value = '''1.000'
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '''4.823'
$ws.Range("E23").Value = '  +5.40%  '
$ws.Range("E24").Value = '  +7.50%  '
$ws.Range("D25").Value = '''9.262'
$ws.Range("D26").Value = '''147.48'
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").Value = '''134.78'
$ws.Range("E27").Value = '  +4.80%  '
$ws.Range("D28").Value = '''17.05'
$ws.Range("E28").Value = '  +4.12%  '
$ws.Range("D29").Value = '''1.993'
$ws.Range("E29").Value = '  +5.92%  '
$ws.Range("D30").Value = '''1.405'
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("D32").Value = '''0.08806'
$ws.Range("E32").Value = '  +5.52%  '
$ws.Range("D33").Value = '''3.997'
$ws.Range("E33").Value = '  +5.08%  '
$ws.Range("D34").Value = '''0.05049'
$ws.Range("E34").Value = '  +2.13%  '
$ws.Range("D35").Value = '''1.138'
$ws.Range("E35").Value = '  +4.52%  '
$ws.Range("D36").Value = '''0.7192'
$ws.Range("E36").Value = '  +7.49%  '
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("D38").Value = '''2.813'
$ws.Range("D39").Value = '''2.268'
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").Value = '''0.9643'
$ws.Range("E40").Value = '  +1.00%  '
$ws.Range("D41").Value = '''0.01692'
$ws.Range("E41").Value = '  +6.45%  '
$ws.Range("D42").Value = '''6.127'
$ws.Range("E42").Value = '  +0.59%  '
$ws.Range("D43").Value = '''0.4303'
$ws.Range("E43").Value = '  +5.54%  '
$ws.Range("D44").Value = '''104.61'
$ws.Range("E44").Value = '  +4.93%  '
$ws.Range("D45").Value = '''0.9989'
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("D46").Value = '''7.606'
$ws.Range("E46").Value = '  +6.48%  '
$ws.Range("E47").Value = '  +5.10%  '
$ws.Range("E48").Value = '  +4.18%  '
$ws.Range("D49").Value = '''33.36'
$ws.Range("E49").Value = '  +5.69%  '
$ws.Range("D50").Value = '''8.391'
$ws.Range("E50").Value = '  +3.58%  '
$ws.Range("D51").Value = '''0.3808'
$ws.Range("E51").Value = '  +5.07%  '
